$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range('D2').Value = '26.294.47'
$ws.Range('E2').Value = '  +0.99%  '
$ws.Range('D3').Value = '1.679.35'
$ws.Range('E3').Value = '  +0.65%  '
$ws.Range('E4').Value = '  +0.27%  '
$ws.Range('D5').NumberFormat = '@'
$ws.Range('D5').Value = '217.88'
$ws.Range('D5').NumberFormat = 'General'
$ws.Range('E5').Value = '  +0.56%  '
$ws.Range('D6').NumberFormat = '@'
$ws.Range('D6').Value = '0.5342'
$ws.Range('D6').NumberFormat = 'General'
$ws.Range('E6').Value = '  +4.82%  '
$ws.Range('E7').Value = '  +0.23%  '
$ws.Range('E8').Value = '  +1.15%  '
$ws.Range('D9').NumberFormat = '@'
$ws.Range('D9').Value = '0.06473'
$ws.Range('D9').NumberFormat = 'General'
$ws.Range('E9').Value = '  +1.08%  '
$ws.Range('D10').NumberFormat = '@'
$ws.Range('D10').Value = '21.93'
$ws.Range('D10').NumberFormat = 'General'
$ws.Range('E11').Value = '  +1.58%  '
$ws.Range('D12').Value = '1.681.54'
$ws.Range('E12').Value = '  +0.50%  '
$ws.Range('D13').NumberFormat = '@'
$ws.Range('D13').Value = '4.526'
$ws.Range('D13').NumberFormat = 'General'
$ws.Range('E13').Value = '  +0.58%  '
$ws.Range('D14').NumberFormat = '@'
$ws.Range('D14').Value = '0.5779'
$ws.Range('D14').NumberFormat = 'General'
$ws.Range('E14').Value = '  -1.09%  '
$ws.Range('D15').NumberFormat = '@'
$ws.Range('D15').Value = '0.000008445'
$ws.Range('D15').NumberFormat = 'General'
$ws.Range('E15').Value = '  -1.13%  '
$ws.Range('D16').NumberFormat = '@'
$ws.Range('D16').Value = '64.78'
$ws.Range('D16').NumberFormat = 'General'
$ws.Range('E16').Value = '  +0.52%  '
$ws.Range('D17').Value = '26.319.22'
$ws.Range('E18').Value = '  -0.81%  '
$ws.Range('E19').Value = '  +0.27%  '
$ws.Range('E20').Value = '  +0.82%  '
$ws.Range('D21').NumberFormat = '@'
$ws.Range('D21').Value = '191.06'
$ws.Range('D21').NumberFormat = 'General'
$ws.Range('D22').NumberFormat = '@'
$ws.Range('D22').Value = '6.210'
$ws.Range('D22').NumberFormat = 'General'
$ws.Range('E22').Value = '  -0.15%  '
$ws.Range('E23').Value = '  +0.21%  '
$ws.Range('D24').NumberFormat = '@'
$ws.Range('D24').Value = '146.15'
$ws.Range('D24').NumberFormat = 'General'
$ws.Range('E24').Value = '  +0.85%  '
$ws.Range('D25').NumberFormat = '@'
$ws.Range('D25').Value = '0.1280'
$ws.Range('D25').NumberFormat = 'General'
$ws.Range('E25').Value = '  +7.13%  '
$ws.Range('D26').NumberFormat = '@'
$ws.Range('D26').Value = '7.836'
$ws.Range('D26').NumberFormat = 'General'
$ws.Range('E26').Value = '  +2.90%  '
$ws.Range('E27').Value = '  +0.12%  '
$ws.Range('D28').NumberFormat = '@'
$ws.Range('D28').Value = '0.06499'
$ws.Range('D28').NumberFormat = 'General'
$ws.Range('E28').Value = '  +1.41%  '
$ws.Range('D29').NumberFormat = '@'
$ws.Range('D29').Value = '1.392'
$ws.Range('D29').NumberFormat = 'General'
$ws.Range('E29').Value = '  +4.26%  '
$ws.Range('D31').NumberFormat = '@'
$ws.Range('D31').Value = '3.579'
$ws.Range('D31').NumberFormat = 'General'
$ws.Range('E31').Value = '  +0.92%  '
$ws.Range('D32').NumberFormat = '@'
$ws.Range('D32').Value = '3.575'
$ws.Range('D32').NumberFormat = 'General'
$ws.Range('E32').Value = '  +1.56%  '
$ws.Range('D33').NumberFormat = '@'
$ws.Range('D33').Value = '1.666'
$ws.Range('D33').NumberFormat = 'General'
$ws.Range('E33').Value = '  +1.18%  '
$ws.Range('E34').Value = '  +1.37%  '
$ws.Range('D35').NumberFormat = '@'
$ws.Range('D35').Value = '0.6163'
$ws.Range('D35').NumberFormat = 'General'
$ws.Range('E35').Value = '  +1.02%  '
$ws.Range('E36').Value = '  +1.38%  '
$ws.Range('D37').NumberFormat = '@'
$ws.Range('D37').Value = '2.712'
$ws.Range('D37').NumberFormat = 'General'
$ws.Range('E37').Value = '  +0.17%  '
$ws.Range('D38').NumberFormat = '@'
$ws.Range('D38').Value = '6.243'
$ws.Range('D38').NumberFormat = 'General'
$ws.Range('E38').Value = '  -0.14%  '
$ws.Range('D39').Value = '1.110.18'
$ws.Range('E39').Value = '  +2.00%  '
$ws.Range('D40').NumberFormat = '@'
$ws.Range('D40').Value = '0.01620'
$ws.Range('D40').NumberFormat = 'General'
$ws.Range('E40').Value = '  +1.11%  '
$ws.Range('D41').NumberFormat = '@'
$ws.Range('D41').Value = '0.8692'
$ws.Range('D41').NumberFormat = 'General'
$ws.Range('E41').Value = '  +0.99%  '
$ws.Range('E42').Value = '  +0.60%  '
$ws.Range('D43').NumberFormat = '@'
$ws.Range('D43').Value = '100.54'
$ws.Range('D43').NumberFormat = 'General'
$ws.Range('E43').Value = '  -0.04%  '
$ws.Range('D44').Value = '1.830.04'
$ws.Range('E44').Value = '  +0.78%  '
$ws.Range('E45').Value = '  -5.14%  '
$ws.Range('D46').NumberFormat = '@'
$ws.Range('D46').Value = '57.09'
$ws.Range('D46').NumberFormat = 'General'
$ws.Range('E46').Value = '  +1.26%  '
$ws.Range('D47').NumberFormat = '@'
$ws.Range('D47').Value = '8.143'
$ws.Range('D47').NumberFormat = 'General'
$ws.Range('E47').Value = '  +1.01%  '
$ws.Range('D48').NumberFormat = '@'
$ws.Range('D48').Value = '0.9995'
$ws.Range('D48').NumberFormat = 'General'
$ws.Range('E48').Value = '  -0.53%  '
$ws.Range('D49').NumberFormat = '@'
$ws.Range('D49').Value = '0.05264'
$ws.Range('D49').NumberFormat = 'General'
$ws.Range('E49').Value = '  +0.57%  '
$ws.Range('D50').NumberFormat = '@'
$ws.Range('D50').Value = '0.4292'
$ws.Range('D50').NumberFormat = 'General'
$ws.Range('E50').Value = '  +0.12%  '
$ws.Range('D51').NumberFormat = '@'
$ws.Range('D51').Value = '6.076'
$ws.Range('D51').NumberFormat = 'General'
$ws.Range('E51').Value = '  +0.73%  '
